$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 9.0623843370648842
$ws.Range("C2").Value = 3.2811222617082194
$ws.Range("D2").Value = 3.0776763777496114
$ws.Range("E2").Value = 3.864068724345799

# Row 3 values
$ws.Range("B3").Value = 7.0653639578236493
$ws.Range("C3").Value = 15.758628099767833
$ws.Range("D3").Value = 15.865637420340249
$ws.Range("E3").Value = -1.8033215746505107

# Update the selection range to match new selection B1:E3
$ws.Range("B1:E3").Select()
